# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# sheets, incrementing each by 1, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 8313
    $ws.Range("F3").Value = 7737
    $ws.Range("F4").Value = 117
}
